# feat: add 2022-Q3 data
#
# 1) Insert a new worksheet "2022-Q3" before the current "2022-Q2" sheet
#    (position 2) and populate it with the fund-holdings detail table.
# 2) Insert a new top row into the "总计" (totals) summary sheet for the
#    "2022-Q3" quarter, pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: new "2022-Q3" worksheet with fund-holdings detail data
# ---------------------------------------------------------------------

$total = $wb.Worksheets.Item(1)
$q2 = $wb.Worksheets.Item(2)

$newSheet = $wb.Worksheets.Add($q2)
$newSheet.Name = "2022-Q3"

# Header row (B1:H1) - text labels, bold/centered style copied from the
# matching header on the "总计" sheet.
$newSheet.Cells.Item(1,2).Value2 = "基金代码"
$newSheet.Cells.Item(1,3).Value2 = "基金名称"
$newSheet.Cells.Item(1,4).Value2 = "基金规模"
$newSheet.Cells.Item(1,5).Value2 = "股票总仓位"
$newSheet.Cells.Item(1,6).Value2 = "仓位占比"
$newSheet.Cells.Item(1,7).Value2 = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value2 = "仓位排名"

$total.Cells.Item(1,2).Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Data rows. Column A (row index) and H (rank) are numbers; B-G are
# stored as literal text (fund codes, names and plain numeric-looking
# figures), matching how the source data was produced.
$rows = @(
    @{ idx=0; code="516620"; name="国泰中证影视主题ETF";         scale="0.94"; pos="99.07"; pct="3.76"; mv="0.0353"; rank=9 },
    @{ idx=1; code="159855"; name="银华中证影视主题ETF";         scale="0.84"; pos="96.84"; pct="3.69"; mv="0.0310"; rank=9 },
    @{ idx=2; code="004890"; name="中邮健康文娱灵活配置混合";   scale="0.41"; pos="86.15"; pct="3.63"; mv="0.0149"; rank=9 },
    @{ idx=3; code="003397"; name="银华体育文化灵活配置混合";   scale="0.32"; pos="81.07"; pct="3.66"; mv="0.0117"; rank=9 }
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r,1).Value2 = $row.idx
    $newSheet.Cells.Item($r,2).Value2 = "'" + $row.code
    $newSheet.Cells.Item($r,3).Value2 = "'" + $row.name
    $newSheet.Cells.Item($r,4).Value2 = "'" + $row.scale
    $newSheet.Cells.Item($r,5).Value2 = "'" + $row.pos
    $newSheet.Cells.Item($r,6).Value2 = "'" + $row.pct
    $newSheet.Cells.Item($r,7).Value2 = "'" + $row.mv
    $newSheet.Cells.Item($r,8).Value2 = $row.rank
    $r = $r + 1
}

# Column A cells (row index) reuse the same style as the "总计" sheet's
# index column.
$total.Cells.Item(2,1).Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Part 2: shift "总计" summary rows down and insert the 2022-Q3 totals
# ---------------------------------------------------------------------

for ($r = 9; $r -ge 3; $r--) {
    $srcRow = $r - 1
    $bVal = $total.Cells.Item($srcRow, 2).Value2
    $cVal = $total.Cells.Item($srcRow, 3).Value2
    $dVal = $total.Cells.Item($srcRow, 4).Value2
    $total.Cells.Item($r, 1).Value2 = ($r - 2)
    $total.Cells.Item($r, 2).Value2 = $bVal
    $total.Cells.Item($r, 3).Value2 = $cVal
    $total.Cells.Item($r, 4).Value2 = $dVal
}

# Row 9 is brand new - copy the index-column style down from row 8.
$total.Cells.Item(8,1).Copy()
$total.Cells.Item(9,1).PasteSpecial(-4122)

$total.Cells.Item(2,1).Value2 = 0
$total.Cells.Item(2,2).Value2 = "2022-Q3"
$total.Cells.Item(2,3).Value2 = 4
$total.Cells.Item(2,4).Value2 = 0.09

# Restore the original active sheet (总计).
$total.Activate()
